$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the rows that no longer belong in the refactored table ---
# Old rows 10-14 were FACASIA, FACASIB, NDAASI, NDBASI, NOTADB.
# Deleting row 10 five times in a row removes all five, shifting
# NCA/NCB/NCECA (old rows 15-17) up to rows 10-12.
$ws.Rows.Item(10).Delete()
$ws.Rows.Item(10).Delete()
$ws.Rows.Item(10).Delete()
$ws.Rows.Item(10).Delete()
$ws.Rows.Item(10).Delete()

$lo = $ws.ListObjects.Item(1)

# --- Add a new row for the "recibo" entry (RECX2) ---
$lo.ListRows.Add() | Out-Null

# --- Copy the formatting of the last "notacredito" row (12) down onto
#     the brand new row 13 so the new row matches the surrounding style ---
$ws.Range("A12:B12").Copy()
$ws.Range("A13:B13").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Add the new "tipo2" column to the table ---
$newCol = $lo.ListColumns.Add()
$newCol.Range.Item(1).Value = "tipo2"

# --- Fix up column B values: the former "notadebito" rows become
#     "factura", and the "notacredito" rows become "test" ---
$ws.Range("B6").Value = "factura"
$ws.Range("B7").Value = "factura"
$ws.Range("B8").Value = "factura"
$ws.Range("B9").Value = "factura"
$ws.Range("B10").Value = "test"
$ws.Range("B11").Value = "test"
$ws.Range("B12").Value = "test"

# --- New row 13 data (RECX2 / test / recibo) ---
$ws.Range("A13").Value = "RECX2"
$ws.Range("B13").Value = "test"

# --- Populate the new column C (tipo2) for every data row ---
$ws.Range("C2").Value = "factura"
$ws.Range("C3").Value = "factura"
$ws.Range("C4").Value = "factura"
$ws.Range("C5").Value = "factura"
$ws.Range("C6").Value = "notadebito"
$ws.Range("C7").Value = "notadebito"
$ws.Range("C8").Value = "notadebito"
$ws.Range("C9").Value = "notadebito"
$ws.Range("C10").Value = "notacredito"
$ws.Range("C11").Value = "notacredito"
$ws.Range("C12").Value = "notacredito"
$ws.Range("C13").Value = "recibo"

# --- Match the header/data cell formatting used by the rest of the table ---
$ws.Range("C1").HorizontalAlignment = -4108
$ws.Range("C1").WrapText = $true
$ws.Range("C2:C13").HorizontalAlignment = -4108

# --- Sheet/column cosmetics to mirror the refactor ---
$ws.Columns.Item(3).ColumnWidth = 12.25

# Best-effort: scroll the view so row 5 is at the top (topLeftCell="A5"
# in the target); keep the existing C13 selection untouched like the
# original file so we don't regress that if scrolling isn't persisted.
$win = $wb.Windows.Item(1)
$win.ScrollRow = 5
$win.ScrollColumn = 1
